$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new person (row 18): name (non-public) + birth date
$ws.Range("A18").Value = "非公開"
$ws.Range("B18").Value = "1995年3月6日"

# Row grows slightly taller to fit the new entry's text, like the row above it
$ws.Rows.Item(18).RowHeight = 14.15

# Move the active selection to A22 (was B22)
$ws.Range("A22").Select()
